$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 68.93939393939394
$ws.Range("C3").Value = 73.48484848484848
$ws.Range("C4").Value = 69.6969696969697
$ws.Range("C5").Value = 74.24242424242425
